# Apply daily-scrape update to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2 (opportunity changed from 1328547 to 1328552) ---
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "1328552"
$ws.Range("A2").Style = "Normal"
$ws.Range("B2").Value = "https://aiesec.org/opportunity/global-talent/1328552"
$ws.Range("C2").Value = "Mobile applicatio"
$ws.Range("F2").Value = "0 applicants"

# --- Update row 3 (opportunity changed from 1328541 to 1328548) ---
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "1328548"
$ws.Range("A3").Style = "Normal"
$ws.Range("B3").Value = "https://aiesec.org/opportunity/global-talent/1328548"
$ws.Range("C3").Value = "Sales"

# --- Remove rows 4-6 entirely (no longer scraped) ---
$ws.Range("A4:H6").EntireRow.Delete()

# --- Adjust column widths to match refreshed layout ---
# (offset by -5/6 to compensate for the ColumnWidth -> stored "width" padding
# so the saved XML ends up with the exact integer widths 20 / 70 / 15 / 16)
$ws.Columns.Item(3).ColumnWidth = 19.1666666666667
$ws.Columns.Item(4).ColumnWidth = 69.1666666666667
$ws.Columns.Item(7).ColumnWidth = 14.1666666666667
$ws.Columns.Item(8).ColumnWidth = 15.1666666666667
